$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A60").Value = "2024-09-29 00:00:00"
$ws.Range("B60").Value = 75550
$ws.Range("C60").Value = 10762.57
$ws.Range("D60").Value = 9524.4
$ws.Range("E60").Value = 7.0117
